$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update customer data (row 2) per new "ModificarCliente" transaction
$ws.Range("C2").Value = "'72145803"
$ws.Range("D2").Value = "'Laynes"
$ws.Range("G2").Value = "'Luis Laynes Castro"
$ws.Range("M2").Value = "3 jul. 2023, 14:55:25"

# Update active selection
$ws.Range("I10").Select()

# Column A widened (best-fit) after the data refresh
$ws.Columns.Item(1).ColumnWidth = 16
